$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving plain decimal values must be pre-formatted as Text so
# Excel stores the literal string (preserving trailing zeros etc.) instead
# of auto-converting the input into a number.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '26.932.21'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '1.550.77'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.57%  '
$ws.Range('D5').Value = '206.99'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '0.485'
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('D8').Value = '0.248'
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('D9').Value = '21.58'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').Value = '0.0587'
$ws.Range('D11').Value = '0.0860'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').Value = '1.773.94'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = '1.552.83'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').Value = '26.938.82'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '61.82'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '215.10'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').Value = '0.0₃0687'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '7.23'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').Value = '9.14'
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').Value = '152.81'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  +2.92%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').Value = '0.103'
$ws.Range('E29').Value = '  +1.33%  '
$ws.Range('E30').Value = '  +0.81%  '
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').Value = '3.20'
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('D33').Value = '1.403.32'
$ws.Range('E33').Value = '  +4.19%  '
$ws.Range('E34').Value = '  +2.69%  '
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').Value = '  +3.05%  '
$ws.Range('D36').Value = '0.956'
$ws.Range('E36').Value = '  +2.27%  '
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').Value = '0.521'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').Value = '0.808'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('D42').Value = '0.990'
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '5.51'
$ws.Range('E43').Value = '  -4.77%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '2.26'
$ws.Range('E44').Value = '  +2.99%  '
$ws.Range('D45').Value = '63.64'
$ws.Range('E45').Value = '  +1.68%  '
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = '1.687.58'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').Value = '86.18'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D50').Value = '0.0952'
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('E51').Value = '  +0.68%  '
